$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new player row (5th data row) with just a name - no URL yet.
$ws.Range("A5").Value = "Anthony Edwards"

# Re-fit column A now that it holds a longer name than before.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where Excel would land after typing into A5 and
# pressing Enter/Down.
$ws.Range("A6").Select() | Out-Null
